$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Delete the old "Joystick shield for adruino" row (row 18) first (highest row number first)
$ws.Rows.Item(18).Delete()

# 2. Delete the old "Uno R3" row (row 11)
$ws.Rows.Item(11).Delete()

# 3. Insert a new blank row at row 5 (shifts everything from old row 5 down to row 6)
$ws.Rows.Item(5).Insert()

# 4. Fill new row 5 with the Arduino Uno R3 entry
$ws.Range("A5").Value = "Arduino Uno R3"
$ws.Range("B5").Value = "A000066"
$ws.Range("C5").Value = "ARDUINO CC"
$ws.Range("D5").Value = "A000066 "
$ws.Range("E5").Value = "core-electronics.com.au"
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 35.85
$ws.Range("H5").Formula = "=F5*G5"
$ws.Range("I5").Value = "Controls the Motors"

Write-Output "done"
